$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicate-ID values in column B (rows 2-4)
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 101
$ws.Range("B4").Value = 3

# Row 5 no longer holds data - clear A5 (keeps formatting) and B5 (removes cell)
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()

# Move active selection to C7
$ws.Range("C7").Select()
